# ProjectIoTGarden_ExcelSheet.xlsx — "Power" sheet: add a 32V->5V and a
# 5V->5V voltage-divider calculation block (mirroring the existing
# 9V->3.3V / 9V->5V blocks), plus a small standalone resistor-sum scratch
# calc (U12:U13, U18).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power")

# ---------------------------------------------------------------------
# Header row (row 4): titles for the two new blocks, styled like the
# existing D4 / J4 headers (yellow fill, centered, bordered).
# ---------------------------------------------------------------------
$ws.Range("O4").Value = "32 volt to 5 volt voltage divider"
$ws.Range("W4").Value = "5 volt to 5 volt voltage divider"

$ws.Range("D4").Copy() | Out-Null
$ws.Range("O4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D4").Copy() | Out-Null
$ws.Range("W4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 32 volt -> 5 volt divider (columns O:Q), mirrors D:F
# ---------------------------------------------------------------------
$ws.Range("O6").Value = "Voltage Input"
$ws.Range("P6").Value = 32
$ws.Range("Q6").Value = "volts"

$ws.Range("O7").Value = "Desired Output"
$ws.Range("P7").Value = 5
$ws.Range("Q7").Value = "Volts"

$ws.Range("O8").Value = "Resistor 2"
$ws.Range("P8").Value = 1000
$ws.Range("Q8").Value = "ohm"

$ws.Range("P9").Value = " "

$ws.Range("O10").Value = "Vout / Vin"
$ws.Range("P10").Formula = "=P7/P6"

$ws.Range("O11").Value = "1/(R1+R2)"
$ws.Range("P11").Formula = "=P10/P8"

$ws.Range("O12").Value = "R1+R2"
$ws.Range("P12").Formula = "=1/P11"

$ws.Range("O14").Value = "Resistor 1"
$ws.Range("P14").Formula = "=P12-P8"
$ws.Range("Q14").Value = "ohms"

$ws.Range("O15").Value = "Current At Output"
$ws.Range("P15").Formula = "=P7/P14"
$ws.Range("Q15").Value = "amps"

# ---------------------------------------------------------------------
# 5 volt -> 5 volt divider (columns W:Y), mirrors D:F / O:Q
# ---------------------------------------------------------------------
$ws.Range("W6").Value = "Voltage Input"
$ws.Range("X6").Value = 5
$ws.Range("Y6").Value = "volts"

$ws.Range("W7").Value = "Desired Output"
$ws.Range("X7").Value = 2
$ws.Range("Y7").Value = "Volts"

$ws.Range("W8").Value = "Resistor 2"
$ws.Range("X8").Value = 100
$ws.Range("Y8").Value = "ohm"

$ws.Range("X9").Value = " "

$ws.Range("W10").Value = "Vout / Vin"
$ws.Range("X10").Formula = "=X7/X6"

$ws.Range("W11").Value = "1/(R1+R2)"
$ws.Range("X11").Formula = "=X10/X8"

$ws.Range("W12").Value = "R1+R2"
$ws.Range("X12").Formula = "=1/X11"

$ws.Range("W14").Value = "Resistor 1"
$ws.Range("X14").Formula = "=X12-X8"
$ws.Range("Y14").Value = "ohms"

$ws.Range("W15").Value = "Current At Output"
$ws.Range("X15").Formula = "=X7/X14"
$ws.Range("Y15").Value = "amps"

# ---------------------------------------------------------------------
# Small standalone scratch calc in column U
# ---------------------------------------------------------------------
$ws.Range("U12").Formula = "=1502"
$ws.Range("U13").Value = 3900
$ws.Range("U18").Formula = "=SUM(U12:U16)"

# ---------------------------------------------------------------------
# Column widths for the new label columns (O, W) — sized to fit their
# header text, like the existing best-fit D / J columns.
# ---------------------------------------------------------------------
$ws.Columns.Item(15).ColumnWidth = 25.75
$ws.Columns.Item(23).ColumnWidth = 26.75

# ---------------------------------------------------------------------
# View state: zoom out a bit and move the selection, as when the user
# scrolled right to review the new blocks.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 52
$ws.Range("AG15").Select() | Out-Null
